$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.681.59"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.287.39"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'506.36"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'129.07"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "2.311.70"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("E13").Value = "  +4.40%  "
$ws.Range("D14").Value = "'23.60"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "2.688.97"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "54.696.04"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "2.251.87"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "'10.62"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'4.19"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'6.64"
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("D22").Value = "'308.16"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "'60.40"
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "'0.993"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'0.151"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "'7.49"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("D28").Value = "'171.46"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "0.0₃0705"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  +5.30%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'18.03"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "'0.909"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").Value = "'36.66"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'132.38"
$ws.Range("E42").Value = "  +5.93%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.42"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'4.86"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "'252.64"
$ws.Range("D46").Value = "'0.0501"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").Value = "'0.0913"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "'0.376"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("E51").Value = "  +0.40%  "
